$d = $word.ActiveDocument

$found0 = $d.Content.Find.Execute("המאמר היומי של מייק - 01.01.25:", $true, $false, $false, $false, $false, $true, 1, $false, "המאמר היומי של מייק - 30.12.24:", 2)
if (-not $found0) { Write-Output "REPLACE 0 FAILED" }
$found1 = $d.Content.Find.Execute("Inference-Aware Fine-Tuning for Best-of-N Sampling in Large Language Models", $true, $false, $false, $false, $false, $true, 1, $false, "Graph Diffusion Policy Optimization", 2)
if (-not $found1) { Write-Output "REPLACE 1 FAILED" }
$found2 = $d.Content.Find.Execute("מתחילים את השנה החדשה עם סקירה של מאמר די מעניין שמציע שיטה לשיפור אימון של מודלי שפה. היתרון הגדול של השיטה היא מאפשרת להתאים את האימון לאופן ההיסק (אינפרנס) ודי ברור שאם אכן עושים זאת בהצלחה זה אמור להניב איכות ההיסק. כלומר אם אנו משתמשים בגישה מסוימת במהלך האינפרנס: למשל לבחור את ״התשובה הטובה ביותר״ מבין N תשובות המודל (המאמר מפתח שיטות רק לגישה זו וקורא לה BoN) או תיקון עצמי (self-correction) אז כדאי לנו להתאים את האינפרנס לכך.", $true, $false, $false, $false, $false, $true, 1, $false, "לפני יומיים סקרתי מאמר על מודלי דיפוזיה המאומנים באמצעות שיטות מעולם למידה עם חיזוקים או RL, אתמול סקרתי מאמר על רשתות נוירונים על גרפים והיום החלטתי לסקור מאמר שמאחד את 3 הדברים האלו (כמעט). המאמר המסוקר היום מציע שיטה לאימון מודל המגנרט גרפים באמצעות מודלי דיפוזיה המאומנים עם שיטות RL (נכון אין כאן GNN בצורתם הטהורה אבל לפחות יש גרפים…", 2)
if (-not $found2) { Write-Output "REPLACE 2 FAILED" }
$found3 = $d.Content.Find.Execute("קודם כל המאמר מנסח שתי פונקצית יעד לאימון inference-aware או IA בקצרה, אחת ל SFT וגם ל-RLHF, שקיבלו שמות IA-SFT ו-IA-RL בהתאמה. עבור IA-SFT אנו ממקסמים את נראות של תשובות המומחים לשאלות מהדאטהסט שיש לנו בהינתן פוליסי האינפרנס I (שזה למעשה BoN). למעשה מאפטמים את הפוליסי (שזה מנגנון חיזוי של LLM או בפשטות LLM עצמו) כדי לעשות את BoN בצורה הטובה ביותר על הדאטהסט שיש לנו. עבור IA-RL המטרה היא לאפטם את הפוליסי (שזה LLM כאמור) תחת טכניקה של אינפרנס I (כלומר BoN) כך שהיא ימקסם את פונקצית תגמול R. ", $true, $false, $false, $false, $false, $true, 1, $false, "קודם כל אנו צריכים להבין איך ניתן למנף מודלי דיפוזיה לגנרוט גרפים. האמת זה די פשוט ודומה לגנרוט תמונות. אתם זוכרים מודלי דיפוזיה מאומנים לגנרט תמונה מרעש טהור (בד״כ) על ידי הורדה הדרגתית של הקומפוננטה הרועשת שלו עד להפיכתו לפיסת דאטה המפלגות לפי ההתפלגות של דאטהסט אימון. זה ממש בגדול ויש גישות חדשות יותר שעושות את זה טיפה אחרת למשל כמו Consistency Models שדיברנו עליהם באחת הסקירות הקודמות. ", 2)
if (-not $found3) { Write-Output "REPLACE 3 FAILED" }
$found4 = $d.Content.Find.Execute("לאחר מכן המאמר מגדיר באופן מדויק מה זה BoN (נוסחה 1) כאשר המטרה היא למקסם את איכות התשובות של המודל כאשר אנו בוחרים את התשובה לפי מה שנקרא verifier score (= ציון לאיכות התשובה). דרך אגב מתווספים כאן שני פרמטרים נוספים שהם מספר התשובות שמהן בוחרים את התשובה הכי טובה וגם טמפרטורת T של מודל השפה. באופן אינטואיטיבי ככל ש T גבוהה יותר (יותר רנדומליות ויצירתיות התשובות) מספר התשובות N צריך לעלות.", $true, $false, $false, $false, $false, $true, 1, $false, "האם אנחנו יכולים לעשות משהו דומה עם גרפים? מתברר שכן. אנו יכולים להתחיל מלדגום גרף באקראי (כלומר הצמתים והקשתות שלו) ולאמן מודל לשנות את הערכים בצמתים ובקשתות כך שהגרף יהפוך להיות ״דומה״ לאחד הגרפים מדאטהסט האימון וגם יקבל ערך גבוה לפי איזה פונקציית תגמול(המאמר גם על RL, זוכרים). ד״א, יש כאן הנחה סמויה שצומת יכול לקבל מספר סופי של ערכים (נגיד מ 0 עד a) וכל קשת יכולה להיות מכמה סוגים (כלומר מ- 0 עד b). כלומר ההתפלגויות שאנו דוגמים מהם הם קטגוריאליות וזה שונה ממה שאנו רגילים לראות במודלי דיפוזיה גנרטיביים עבור התמונות.", 2)
if (-not $found4) { Write-Output "REPLACE 4 FAILED" }
$found5 = $d.Content.Find.Execute(" כאן יש לנו כאן הטרייד-אוף הקלאסי בין exploration ל-exploitation. ככל ש- T גבוה יותר אנו מבצעים יותר exploration (תשובות מגוונות יותר) ואילו T קטן יותר מאפשר לנו ״להנות״ ממה שלמדנו עד עכשיו (בחירה של N משפיע על הטרייד-אוף באופן הפוך מ-T).", $true, $false, $false, $false, $false, $true, 1, $false, "כמובן מיד עולות כמה שאלות בנוגע לתהליך הזה?", 2)
if (-not $found5) { Write-Output "REPLACE 5 FAILED" }
$found6 = $d.Content.Find.Execute("אוקיי, אבל עדיין בבעיית אופטימיזציה של IA-SFT יש לנו את argmax (משתמשים בו לבחירה של התשובה הטובה ביותר) וזה מאוד מקשה על פתרונה למרות שיש לנו שיטות שערוך argmax באמצעות softmax וגם Gumbel softmax עדיין שיטות אלו אינן מדויקות וגם כבדות חישובית (לטענת המאמר). אז המחברים משתמשים בטריק מאוד מפורסם ב-ML - קירוב של פונקציית יעד עם קירוב וריאציוני שהופך אותה (את הלוג שלה) לסכום של הפוליסי (הסתברות של תשובה y בהינתן שאלה x עם המודל) ושל איבר רגולריזציה הנקרא inference aware. איבר זה הוא למעשה הוא win-rate של תשובה y על שאלה x על המודל הנוכחי כאשר הערך של כל זוג (x, y) מחושב עם verifier score r (עם קבוע נרמול).", $true, $false, $false, $false, $false, $true, 1, $false, "איך דוגמים גרף באקראי במהלך האינפרנס (זה נושא עתיק ונחקר רבות עלי ידי מתמטיקאים ובפרט על ידי ארדוש, המאמר לא מתעמק בזה יותר מדי). דרך אגב במהלך האימון אנו לוקחים גרף מהדאטהסט ומרעישים אותו עלי ידי ״שינוים אקראיים״ בערכי הצמתים ובסוגי הקשתות", 2)
if (-not $found6) { Write-Output "REPLACE 6 FAILED" }
$found7 = $d.Content.Find.Execute("ב- IA-RL הסיפור קצת מסתבך והמחברים משתמשים בתוצאה מאחד המאמרים של ג'ו שולמן (cto לשעבר של openai) עם סרגיי לווין ופיטר אבל האגדיים כדי לקבל שערוך לגרדיאנט שמקבל צורה דומה לאלגוריתם הישן והידוע REINFORCE כלומר המכפלה של הלוג של הפוליסי עם פונקציית תגמול (״ממורכזת״ עם התוחלת של פונקציית תגמול להקטנת השונות). המאמר גם דן במקרים מעניינים של אופטימיזציה של פונקצית היעד של IA-RL לכמה צורות של verifier score r (למשל בינארי).", $true, $false, $false, $false, $false, $true, 1, $false, "איך משווים גרפים, כלומר איך מבינים שגרף שקיבלנו במהלך הגנרוט הוא דומה לגרף מהדאטהסט? יש מספר רב גישות להשוות גרפים על ידי השוואה של התת-גרפים שלהם או להשוות את הלפלסיאן שלהם למשל.", 2)
if (-not $found7) { Write-Output "REPLACE 7 FAILED" }
$found8 = $d.Content.Find.Execute("מאמר די כבד מתמטית ניסיתי (לפי מיטב יכולתי) להנגיש לכם אותו טיפה…", $true, $false, $false, $false, $false, $true, 1, $false, "בחירה של פונקצית reward בדומיין הגרפים לא טריוויאלית בכלל. למשל למשימות גנרוט גרפים למולקולות חדשות אחד המדדים לאיכות הגרף המגונרט הוא חדשנותו יחסית לדברים הקיימים, יעילותו בטיפול במחלה מסוימות או פיזיביליות של סינטוזו (synthetic accessibility). ניתן לבחור reward גם בתור פונקצית דמיון לגרפים הקיימים.", 2)
if (-not $found8) { Write-Output "REPLACE 8 FAILED" }
$found9 = $d.Content.Find.Execute("https://arxiv.org/abs/2412.15287", $true, $false, $false, $false, $false, $true, 1, $false, "אוקיי, אז יש לנו פונקציה להשוואת הגרפים C ופונקצית תגמול לשערוך איכות הגרף r - איך אנו מאמנים מודל דיפוזיה. האמת בצורה די דומה לזו שתיארתי בסקירת של לפני 3 ימים של המאמר: RL for Consistency Models: Faster Reward Guided Text-to-Image Generation. ", 2)
if (-not $found9) { Write-Output "REPLACE 9 FAILED" }

# Append new paragraphs after the last paragraph in the document
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "קודם כל אנו צריכים להגדיר את Markov Decision Process עבור אימון מודל דיפוזיה על גרפים. ומתברר שהוא ממש דומה למאמר שהזכרתי:"
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "המצב s_t בתור זוג של גרף מגונרטת באיטרציה T-t וגם ערך T-t" + [char]11 + "הפעולה a_t היא הגרף באיטרציה T−t−1"
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "הפוליסי (הסתברות של a_t בהינתן s_t) היא זו פונקצית התפלגות מותנית של גרף מאיטרציה T−t−1 בהינתן גרף באיטרציה T-t"
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "המצב ההתחלתי הוא גרף אקראי באיטרציה T ופונקציית תגמול r שנתונה לנו המחושבת על הגרף הסופי באיטרציה 0"
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "המאמר מציע שתי שיטות לאימון של מודל דיפוזיה לגנרוט גרפים: הראשונה היא REINFORCE הקלאסי שהיא למעשה שיטת policy gradient הממקסת פוליסים בעלי תגמול גבוה. מעשה אנו דוגמים K איטרציה בין 1 ל T וממקסמים מכפלה ממוצעת (על K דגימות) של פונקציית הפוליסו (ונקצית התפלגות מותנית של גרף מאיטרציה T−t−1 בהינתן גרף באיטרציה T-t) והתגמול עבור הגרף המגונרט (באיטרציה 0)."
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "השיטה השנייה המוצעת היא Policy Optimization כאשר במקום למקסם את הפוליסי בצורתו הטהורה אנו ממקסמים הסתברות גנרוט גרף G_0 מהדאטהסט (שאותו מרעישים והמודל ״מסיר״ ממנו את הרעש) מוכפלת בתגמול עבור הגרף הנוצר. גם כאן יש מיצוע על K איטרציות שמהם נבנה שערוך של G_0."
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "זהו זה - סקירה קצת כבדה, מקווה שהצלחתם להבין משהו ממנה…"
$tail.InsertParagraphAfter()
$lastIndex = $d.Paragraphs.Count
$tail = $d.Paragraphs($lastIndex).Range
$tail.Text = "https://arxiv.org/abs/2402.16302"

Write-Output "Done. Paragraphs count: $($d.Paragraphs.Count)"
